$d = $word.ActiveDocument

# 1. Update the date text.
$d.Content.Find.Execute("March 29, 2014", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "May 21, 2014", 2)

# 2. Move the "_GoBack" bookmark from its old location (after "which" in the
#    body paragraph) to the end of the date line (after the new date text).
$d.Bookmarks.Item("_GoBack").Delete()

# Find the end of the date paragraph's text.
$dateRange = $d.Content
$dateRange.Find.Execute("May 21, 2014", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$endPos = $dateRange.End

# A collapsed range sitting exactly at a paragraph's last text gap (just
# before its paragraph mark) gets mis-anchored by Bookmarks.Add, so insert a
# throwaway placeholder character first to make that insertion point
# interior, add the bookmark there, then remove the placeholder again.
$insPoint = $d.Range($endPos, $endPos)
$insPoint.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($endPos, $endPos + 1)
$placeholder.Delete()
